# "Generate Report for Handback"
# Refreshes the timestamps / priority recorded for the 9a280f3c and
# 9bd5e22b handback rows across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" (column G)
# rows 4 (9a280f3c...) and 5 (9bd5e22b...)
$overview.Range("G4").Value = "2016-08-31 20:19:52"
$overview.Range("G5").Value = "2016-08-31 20:19:52"

# zh-cn sheet, rows 4 and 5 (9a280f3c... / 9bd5e22b...)
# Priority (column E): ht -> mt
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("E5").Value = "mt"

# Correspond Handoff Datetime (column H)
$zhcn.Range("H4").Value = "2016-08-31 20:19:47"
$zhcn.Range("H5").Value = "2016-08-31 20:19:47"

# Correspond Handback DateTime (column K)
$zhcn.Range("K4").Value = "2016-08-31 20:20:21"
$zhcn.Range("K5").Value = "2016-08-31 20:20:21"

# de-de sheet, rows 4 and 5 (9a280f3c... / 9bd5e22b...)
# Correspond Handback DateTime (column K)
$dede.Range("K4").Value = "2016-08-31 20:20:28"
$dede.Range("K5").Value = "2016-08-31 20:20:28"
